$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 91; this shifts the existing rows 91-113 down to 92-114,
# preserving formatting (e.g. the date style on column D) via the row above.
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new record's data.
$ws.Range("A91").Value = 10
$ws.Range("B91").Value = "Vega Modelo de Temuco"
$ws.Range("C91").Value = "La Araucanía"
$ws.Range("D91").Value = 45135
$ws.Range("E91").Value = 9
$ws.Range("F91").Value = "Fruta"
$ws.Range("G91").Value = 100108
$ws.Range("H91").Value = "Tropicales y subtropicales"
$ws.Range("I91").Value = 100108007
$ws.Range("J91").Value = "Coco"
$ws.Range("K91").Value = "Sin especificar"
$ws.Range("L91").Value = "Primera"
$ws.Range("M91").Value = 40
$ws.Range("N91").Value = 30000
$ws.Range("O91").Value = 30000
$ws.Range("P91").Value = 30000
$ws.Range("Q91").Value = "$/malla 20 unidades"
$ws.Range("R91").Value = "Perú"
$ws.Range("S91").Value = 1500
$ws.Range("T91").Value = 20
